# Update "Transaksi Penjualan" data in Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: bayar_tanggal 45760 -> 45759
$ws.Range("F2").Value = 45759

# Row 3: bayar_tanggal 45761 -> 45760
$ws.Range("F3").Value = 45760

# Row 4: penjualan_id 3 -> 9, jumlah_bayar 70000 -> 40000, bayar_tanggal 45762 -> 45760
$ws.Range("A4").Value = 9
$ws.Range("C4").Value = 40000
$ws.Range("F4").Value = 45760

# Update the active selection to C7 (matches last saved cursor position)
$ws.Range("C7").Select()
